$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume cells we are about to rewrite to stay text,
# so strings like "51.847.49" or "  +5.06%  " are not reinterpreted as numbers/dates.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

$ws.Range('D2').Value = '51.847.49'
$ws.Range('E2').Value = '  +5.06%  '
$ws.Range('D3').Value = '2.777.56'
$ws.Range('E3').Value = '  +5.90%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '116.21'
$ws.Range('E5').Value = '  +4.92%  '
$ws.Range('D6').Value = '335.44'
$ws.Range('E6').Value = '  +3.42%  '
$ws.Range('D7').Value = '0.539'
$ws.Range('E7').Value = '  +3.28%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.574'
$ws.Range('E9').Value = '  +6.31%  '
$ws.Range('D10').Value = '42.09'
$ws.Range('E10').Value = '  +7.45%  '
$ws.Range('D11').Value = '0.0873'
$ws.Range('E11').Value = '  +8.48%  '
$ws.Range('D12').Value = '20.07'
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('E13').Value = '  +2.43%  '
$ws.Range('D14').Value = '7.61'
$ws.Range('E14').Value = '  +4.52%  '
$ws.Range('D15').Value = '3.216.73'
$ws.Range('E15').Value = '  +5.88%  '
$ws.Range('D16').Value = '2.832.47'
$ws.Range('E16').Value = '  +7.37%  '
$ws.Range('D17').Value = '0.886'
$ws.Range('E17').Value = '  +4.60%  '
$ws.Range('D18').Value = '51.696.93'
$ws.Range('E18').Value = '  +4.81%  '
$ws.Range('D19').Value = '3.25'
$ws.Range('E19').Value = '  +12.34%  '
$ws.Range('D20').Value = '13.37'
$ws.Range('E20').Value = '  +3.82%  '
$ws.Range('D21').Value = '6.94'
$ws.Range('E21').Value = '  +4.30%  '
$ws.Range('D22').Value = '0.0₃0978'
$ws.Range('E22').Value = '  +4.07%  '
$ws.Range('D23').Value = '277.34'
$ws.Range('E23').Value = '  +4.31%  '
$ws.Range('D24').Value = '70.01'
$ws.Range('E24').Value = '  +2.07%  '
$ws.Range('D25').Value = '2.75'
$ws.Range('E25').Value = '  +9.25%  '
$ws.Range('D26').Value = '26.74'
$ws.Range('E26').Value = '  +3.48%  '
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').Value = '10.18'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').Value = '2.22'
$ws.Range('E29').Value = '  +1.01%  '
$ws.Range('D30').Value = '0.142'
$ws.Range('E30').Value = '  +3.88%  '
$ws.Range('D31').Value = '35.09'
$ws.Range('E31').Value = '  +2.35%  '
$ws.Range('D32').Value = '50.11'
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('D33').Value = '5.63'
$ws.Range('E33').Value = '  +3.28%  '
$ws.Range('D34').Value = '0.0819'
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('D35').Value = '2.12'
$ws.Range('E35').Value = '  +5.08%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').Value = '19.02'
$ws.Range('E37').Value = '  +1.46%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '3.29'
$ws.Range('E38').Value = '  +7.05%  '
$ws.Range('D39').Value = '4.95'
$ws.Range('E39').Value = '  +1.25%  '
$ws.Range('D40').Value = '2.70'
$ws.Range('E40').Value = '  +26.27%  '
$ws.Range('D41').Value = '0.0356'
$ws.Range('E41').Value = '  +10.78%  '
$ws.Range('D42').Value = '23.69'
$ws.Range('E42').Value = '  +5.84%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '2.34'
$ws.Range('E43').Value = '  +7.30%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '127.12'
$ws.Range('E44').Value = '  -0.40%  '
$ws.Range('D45').Value = '0.114'
$ws.Range('E45').Value = '  +3.19%  '
$ws.Range('D46').Value = '2.100.26'
$ws.Range('E46').Value = '  +3.13%  '
$ws.Range('D47').Value = '3.30'
$ws.Range('E47').Value = '  +3.80%  '
$ws.Range('D48').Value = '2.23'
$ws.Range('E48').Value = '  +3.51%  '
$ws.Range('D49').Value = '5.54'
$ws.Range('E49').Value = '  +7.16%  '
$ws.Range('D50').Value = '0.911'
$ws.Range('E50').Value = '  +23.26%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').Value = '8.88'
$ws.Range('E51').Value = '  +0.56%  '
